# Append two new match rows (124 and 125) to the Romania Liga-1 2023-2024 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        Row   = 124
        Idx   = 123
        Date  = 45241.73958333334
        Home  = "Poli Iasi"
        HG    = 3
        Away  = "CFR Cluj"
        AG    = 3
        J     = 4.18
        K     = "08/11/2023 08:12"
        L     = 4.56
        M     = "11/11/2023 17:42"
        N     = 3.48
        O     = "08/11/2023 08:12"
        P     = 3.5
        Q     = "11/11/2023 17:36"
        R     = 1.89
        S     = "08/11/2023 08:12"
        T     = 1.84
        U     = "11/11/2023 17:42"
        V     = "https://www.betexplorer.com/football/romania/liga-1/poli-iasi-cfr-cluj/nmmq1uIo/"
    },
    @{
        Row   = 125
        Idx   = 124
        Date  = 45241.83333333334
        Home  = "Farul Constanta"
        HG    = 1
        Away  = "FC Hermannstadt"
        AG    = 1
        J     = 1.81
        K     = "08/11/2023 08:12"
        L     = 1.89
        M     = "11/11/2023 19:59"
        N     = 3.51
        O     = "08/11/2023 08:12"
        P     = 3.53
        Q     = "11/11/2023 19:57"
        R     = 4.53
        S     = "08/11/2023 08:12"
        T     = 4.25
        U     = "11/11/2023 19:59"
        V     = "https://www.betexplorer.com/football/romania/liga-1/farul-constanta-fc-hermannstadt/AV30xduo/"
    }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Idx
    $ws.Cells.Item($row, 1).Style = $ws.Cells.Item($row - 1, 1).Style

    $ws.Cells.Item($row, 2).Value = "romania"
    $ws.Cells.Item($row, 3).Value = "liga-1"
    $ws.Cells.Item($row, 4).Value = "2023-2024"

    $ws.Cells.Item($row, 5).Value = $r.Date
    $ws.Cells.Item($row, 5).Style = $ws.Cells.Item($row - 1, 5).Style

    $ws.Cells.Item($row, 6).Value = $r.Home
    $ws.Cells.Item($row, 7).Value = $r.HG
    $ws.Cells.Item($row, 8).Value = $r.Away
    $ws.Cells.Item($row, 9).Value = $r.AG

    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
}
